$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set DOC_init value (B7) and POC_init value (B8)
$ws.Range("B7").Value = 2.9
$ws.Range("B8").Value = 0.29

# Update the active cell selection shown in the sheet view
$ws.Activate()
$ws.Range("G10").Select()
